$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Price" column cells to store values as literal text (matching the
# original inlineStr cells) instead of letting Excel auto-convert numeric-looking
# strings (e.g. "1.00" -> 1, "0.600" -> 0.6, "68.986.05" stays text anyway).
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D21", "D23", "D24", "D25", "D26", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the refreshed price / 1h-volume-change / coin values cell by cell, in the
# same order as the upstream data refresh.
$ws.Range("D2").Value = '68.986.05'
$ws.Range("E2").Value = '  -4.02%  '
$ws.Range("D3").Value = '3.515.10'
$ws.Range("E3").Value = '  -4.65%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '579.54'
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").Value = '174.17'
$ws.Range("E6").Value = '  -2.72%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.509.46'
$ws.Range("E8").Value = '  -4.63%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E10").Value = '  -6.35%  '
$ws.Range("E11").Value = '  +8.42%  '
$ws.Range("D12").Value = '0.600'
$ws.Range("E12").Value = '  -2.45%  '
$ws.Range("D13").Value = '47.29'
$ws.Range("E13").Value = '  -5.36%  '
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").Value = '670.65'
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").Value = '4.081.21'
$ws.Range("E16").Value = '  -4.67%  '
$ws.Range("D17").Value = '8.81'
$ws.Range("E17").Value = '  -1.82%  '
$ws.Range("D18").Value = '3.532.73'
$ws.Range("E18").Value = '  -3.93%  '
$ws.Range("D19").Value = '68.988.78'
$ws.Range("E19").Value = '  -4.06%  '
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").Value = '17.56'
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D23").Value = '0.906'
$ws.Range("E23").Value = '  -3.70%  '
$ws.Range("D24").Value = '16.29'
$ws.Range("E24").Value = '  -8.64%  '
$ws.Range("D25").Value = '98.35'
$ws.Range("E25").Value = '  -4.99%  '
$ws.Range("D26").Value = '3.87'
$ws.Range("E26").Value = '  -4.28%  '
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  -6.64%  '
$ws.Range("D30").Value = '9.43'
$ws.Range("E30").Value = '  -7.26%  '
$ws.Range("D31").Value = '32.95'
$ws.Range("E31").Value = '  -7.10%  '
$ws.Range("D32").Value = '3.22'
$ws.Range("E32").Value = '  -7.58%  '
$ws.Range("D33").Value = '8.76'
$ws.Range("E33").Value = '  -4.76%  '
$ws.Range("D34").Value = '7.30'
$ws.Range("E34").Value = '  -1.14%  '
$ws.Range("E35").Value = '  -4.55%  '
$ws.Range("D36").Value = '577.82'
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").Value = '10.94'
$ws.Range("E37").Value = '  -3.26%  '
$ws.Range("E38").Value = '  -14.23%  '
$ws.Range("E39").Value = '  -4.09%  '
$ws.Range("D40").Value = '57.08'
$ws.Range("E40").Value = '  -4.31%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D42").Value = '0.337'
$ws.Range("E42").Value = '  -3.39%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0438'
$ws.Range("E43").Value = '  -5.20%  '
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = '0.137'
$ws.Range("E44").Value = '  -5.72%  '
$ws.Range("D45").Value = '3.419.50'
$ws.Range("E45").Value = '  -8.83%  '
$ws.Range("D46").Value = '33.48'
$ws.Range("E46").Value = '  -5.57%  '
$ws.Range("D47").Value = '0.0₃0706'
$ws.Range("E47").Value = '  -8.72%  '
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("E49").Value = '  -6.94%  '
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("D51").Value = '131.21'
$ws.Range("E51").Value = '  -2.22%  '

# Reset the Price-column style back to the sheet default now that the literal
# text is committed, so no stray style index is left referenced on these cells.
foreach ($ref in $priceCells) {
    $ws.Range($ref).Style = "Normal"
}
